# Insert a new weekly price-report row for "Arveja Verde" (Feria Lagunitas de
# Puerto Montt) at row 94, pushing the existing rows 94-157 down to 95-158.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 94..157 down by one row.
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new record.
$ws.Cells.Item(94, 1).Value2 = 4
$ws.Cells.Item(94, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(94, 3).Value2 = "Los Lagos"
$ws.Cells.Item(94, 4).Value2 = 45086
$ws.Cells.Item(94, 5).Value2 = 10
$ws.Cells.Item(94, 6).Value2 = 100112022
$ws.Cells.Item(94, 7).Value2 = "Arveja Verde"
$ws.Cells.Item(94, 8).Value2 = "Perfection"
$ws.Cells.Item(94, 9).Value2 = "Primera"
$ws.Cells.Item(94, 10).Value2 = 40
$ws.Cells.Item(94, 11).Value2 = 43000
$ws.Cells.Item(94, 12).Value2 = 43000
$ws.Cells.Item(94, 13).Value2 = 43000
$ws.Cells.Item(94, 14).Value2 = "`$/malla 25 kilos"
$ws.Cells.Item(94, 15).Value2 = "Provincia de Huasco"
$ws.Cells.Item(94, 16).Value2 = 1720
$ws.Cells.Item(94, 17).Value2 = 25
$ws.Cells.Item(94, 18).Value2 = "Hortaliza"
